# Fruta / hortaliza, semanal
# Insert a new weekly record for Femacal de La Calera - Mango as row 206,
# pushing the existing rows 206:244 down to 207:245.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 206 (shifts rows 206..244 -> 207..245)
$ws.Rows(206).Insert()

# Populate the newly inserted row 206 with the new weekly record
$ws.Cells.Item(206, 1).Value = 3
$ws.Cells.Item(206, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(206, 3).Value = "Coquimbo"
$ws.Cells.Item(206, 4).Value = 44505
$ws.Cells.Item(206, 5).Value = 5
$ws.Cells.Item(206, 6).Value = "Fruta"
$ws.Cells.Item(206, 7).Value = 100108
$ws.Cells.Item(206, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(206, 9).Value = 100108002
$ws.Cells.Item(206, 10).Value = "Mango"
$ws.Cells.Item(206, 11).Value = "Sin especificar"
$ws.Cells.Item(206, 12).Value = "Primera"
$ws.Cells.Item(206, 13).Value = 456
$ws.Cells.Item(206, 14).Value = 7000
$ws.Cells.Item(206, 15).Value = 7000
$ws.Cells.Item(206, 16).Value = 7000
$ws.Cells.Item(206, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(206, 18).Value = "Perú"
$ws.Cells.Item(206, 19).Value = 1750
$ws.Cells.Item(206, 20).Value = 4
